# Fill in the date ("23.09") for the "ЛР02" row of the schedule table and
# shade that date cell with the "Accent 1, Lighter 80%" theme-like fill
# (DBE5F1), matching the formatting already used on the other date cells
# in this table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 9 / Column 1 (1-based, Word Table.Cell indexing) is the empty date
# cell directly above the "ЛР02" / "Приймальні випробування (формування
# плану приймальних випробувань)." row.
$cell = $t.Cell(9, 1)

# Shade the cell (w:fill="DBE5F1" - Accent 1, Lighter 80%).
$cell.Shading.BackgroundPatternColor = 15853019

# Type the date into the (currently empty) cell paragraph.
$cell.Range.InsertBefore("23.09")

# Match the font formatting used throughout the table (Times New Roman,
# 14pt / sz 28).
$cell.Range.Font.Name = "Times New Roman"
$cell.Range.Font.Size = 14
$cell.Range.Font.SizeBi = 14

Write-Output "Set date cell (row 9, col 1) to 23.09 with shading DBE5F1"
